$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.673.37"
$ws.Range("E2").Value = "  -1.44%  "

$ws.Range("D3").Value = "2.074.03"
$ws.Range("E3").Value = "  -2.24%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.49"
$ws.Range("E5").Value = "  -0.68%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.57"
$ws.Range("E7").Value = "  +0.53%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.394"
$ws.Range("E9").Value = "  +0.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0783"
$ws.Range("E10").Value = "  -0.09%  "

$ws.Range("E11").Value = "  +3.20%  "

$ws.Range("D12").Value = "2.378.55"
$ws.Range("E12").Value = "  -2.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.76"
$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.96"
$ws.Range("E14").Value = "  -2.42%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.774"
$ws.Range("E15").Value = "  -1.78%  "

$ws.Range("E16").Value = "  +2.21%  "

$ws.Range("D17").Value = "2.086.88"
$ws.Range("E17").Value = "  -1.73%  "

$ws.Range("D18").Value = "37.638.53"
$ws.Range("E18").Value = "  -1.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.11"
$ws.Range("E19").Value = "  -1.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.41"
$ws.Range("E20").Value = "  +1.16%  "

$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +0.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.40"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("E25").Value = "  -3.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.52"
$ws.Range("E26").Value = "  +1.49%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.06"
$ws.Range("E27").Value = "  +0.58%  "

$ws.Range("E28").Value = "  -2.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.48"
$ws.Range("E29").Value = "  -0.59%  "

$ws.Range("E30").Value = "  -2.93%  "

$ws.Range("E31").Value = "  +1.92%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.68"
$ws.Range("E32").Value = "  +0.53%  "

$ws.Range("E33").Value = "  +0.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.65"

$ws.Range("E35").Value = "  -5.86%  "

$ws.Range("E36").Value = "  -0.56%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.40"
$ws.Range("E37").Value = "  -2.26%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.34"
$ws.Range("E39").Value = "  -2.01%  "

$ws.Range("B40").Value = "Cronos"
$ws.Range("C40").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0974"
$ws.Range("E40").Value = "  -3.04%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.45"
$ws.Range("E41").Value = "  +2.05%  "

$ws.Range("E42").Value = "  -2.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0215"
$ws.Range("E43").Value = "  -0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.72"
$ws.Range("E44").Value = "  +6.20%  "

$ws.Range("D45").Value = "1.435.54"
$ws.Range("E45").Value = "  -2.10%  "

$ws.Range("E46").Value = "  -0.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.21"
$ws.Range("E47").Value = "  +2.18%  "

$ws.Range("E48").Value = "  -0.72%  "

$ws.Range("E49").Value = "  +1.26%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.01"
$ws.Range("E50").Value = "  -1.62%  "

$ws.Range("D51").Value = "2.263.41"
$ws.Range("E51").Value = "  -2.28%  "

Write-Host "Updated cells successfully"